# Update TPM-derived values in the Icam5-Itgb2 LR-pairs sheet.
# This reflects a recalculation of the NATMI output with new TPM values;
# the underlying cluster identities / cell counts (columns A-D, K, L)
# are unchanged, but most of the derived expression / specificity
# columns (E-J, M-T) are recomputed for rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending=ECs, Target=ECs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2577906666666667
$ws.Range("H2").Value = 0.7733719999999999
$ws.Range("I2").Value = 0.1202607703685643
$ws.Range("J2").Value = 0.1202607703685642
$ws.Range("M2").Value = 0.029424
$ws.Range("N2").Value = 0.08827199999999999
$ws.Range("O2").Value = 0.1473063425232919
$ws.Range("P2").Value = 0.1473063425232919
$ws.Range("Q2").Value = 0.007585232575999998
$ws.Range("R2").Value = 0.06826709318399998
$ws.Range("S2").Value = 0.01771517423202668
$ws.Range("T2").Value = 0.01771517423202668

# Row 3 (Sending=ECs, Target=FAPs)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2577906666666667
$ws.Range("H3").Value = 0.7733719999999999
$ws.Range("I3").Value = 0.1202607703685643
$ws.Range("J3").Value = 0.1202607703685642
$ws.Range("O3").Value = 0.852693657476708
$ws.Range("P3").Value = 0.852693657476708
$ws.Range("Q3").Value = 0.04390767971866667
$ws.Range("R3").Value = 0.395169117468
$ws.Range("S3").Value = 0.1025455961365376
$ws.Range("T3").Value = 0.1025455961365375

# Row 4 (Sending=FAPs, Target=ECs)
$ws.Range("I4").Value = 0.8797392296314358
$ws.Range("J4").Value = 0.8797392296314357
$ws.Range("M4").Value = 0.029424
$ws.Range("N4").Value = 0.08827199999999999
$ws.Range("O4").Value = 0.1473063425232919
$ws.Range("P4").Value = 0.1473063425232919
$ws.Range("Q4").Value = 0.05548797535999999
$ws.Range("R4").Value = 0.49939177824
$ws.Range("S4").Value = 0.1295911682912653
$ws.Range("T4").Value = 0.1295911682912652

# Row 5 (Sending=FAPs, Target=FAPs)
$ws.Range("I5").Value = 0.8797392296314358
$ws.Range("J5").Value = 0.8797392296314357
$ws.Range("O5").Value = 0.852693657476708
$ws.Range("P5").Value = 0.852693657476708
$ws.Range("S5").Value = 0.7501480613401705
$ws.Range("T5").Value = 0.7501480613401704
